$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ASV_rank (column G) values for rows 9-12 from 42 to 41
$ws.Range("G9").Value = 41
$ws.Range("G10").Value = 41
$ws.Range("G11").Value = 41
$ws.Range("G12").Value = 41
